$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell in this sheet stores values as text (inlineStr) even when the
# text looks like a plain number (prices, volumes). Plain `Range.Value =
# "123.45"` assignment through this COM layer coerces numeric-looking
# strings to the Number type, which would change the cell's stored type.
# Prefixing the literal with a single-quote forces Excel's "text" entry
# mode (matches typing '123.45 into a cell), and resetting .Style to
# "Normal" afterwards clears the transient quote-prefix style flag so the
# cell's style index is left exactly as it was before the edit.
function Set-TextValue {
    param($range, [string]$value)
    $range.Value = "'" + $value
    $range.Style = "Normal"
}


Set-TextValue $ws.Range("D2") "242.72"
Set-TextValue $ws.Range("D3") "23.44"
Set-TextValue $ws.Range("D4") "5.757"
Set-TextValue $ws.Range("D5") "0.05840"
Set-TextValue $ws.Range("D6") "3.414"
Set-TextValue $ws.Range("D7") "6.496"
Set-TextValue $ws.Range("D8") "1.320"
Set-TextValue $ws.Range("D9") "0.7983"
Set-TextValue $ws.Range("B10") "WazirX"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D10") "0.1476"
Set-TextValue $ws.Range("E10") "9WazirXWRX"
Set-TextValue $ws.Range("B11") "MandalaExchangeToken"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D11") "0.07707"
Set-TextValue $ws.Range("E11") "10MandalaExchangeTokenMDX"
Set-TextValue $ws.Range("B12") "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D12") "0.03258"
Set-TextValue $ws.Range("E12") "11LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue $ws.Range("B13") "BitrueCoin"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.03014"
Set-TextValue $ws.Range("E13") "12BitrueCoinBTR"
Set-TextValue $ws.Range("B14") "BitMartToken"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.09213"
Set-TextValue $ws.Range("E14") "13BitMartTokenBMX"
Set-TextValue $ws.Range("B15") "MCDex"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D15") "3.571"
Set-TextValue $ws.Range("E15") "14MCDexMCB"
Set-TextValue $ws.Range("B16") "BitForexToken"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D16") "0.001665"
Set-TextValue $ws.Range("E16") "15BitForexTokenBF"
Set-TextValue $ws.Range("B17") "CoinExToken"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D17") "0.04757"
Set-TextValue $ws.Range("E17") "16CoinExTokenCET"
Set-TextValue $ws.Range("B18") "One"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D18") "0.0006027"
Set-TextValue $ws.Range("E18") "17OneONE"
Set-TextValue $ws.Range("D19") "0.006254"
Set-TextValue $ws.Range("D20") "0.005481"
Set-TextValue $ws.Range("D21") "0.001070"
Set-TextValue $ws.Range("D23") "3.707"
Set-TextValue $ws.Range("D25") "0.3323"
Set-TextValue $ws.Range("D26") "0.1252"
Set-TextValue $ws.Range("D27") "0.0006277"
Set-TextValue $ws.Range("D40") "0.04329"
Set-TextValue $ws.Range("D41") "0.007044"
Set-TextValue $ws.Range("D42") "0.1054"
Set-TextValue $ws.Range("D43") "0.003410"
Set-TextValue $ws.Range("D44") "0.008682"
Set-TextValue $ws.Range("D45") "0.002463"
Set-TextValue $ws.Range("D46") "0.00005759"
Set-TextValue $ws.Range("D48") "0.7861"
Set-TextValue $ws.Range("D49") "0.1074"
Set-TextValue $ws.Range("E49") "48BOLOBOLOBestin24h"
Set-TextValue $ws.Range("D50") "0.00002102"
Set-TextValue $ws.Range("D51") "0.01011"
